$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.794156074523926
$ws.Range("B1").Value = 2.578701972961426
$ws.Range("C1").Value = 2.864722967147827
$ws.Range("D1").Value = 3.365700244903564
$ws.Range("E1").Value = 3.112615346908569
